{"js": "// The signature block's date line currently reads \"1st October 2024\" and is\n// split across three runs (\"1\", \"st\", \" October 2024\") that all share the\n// SigDate character style. Replace that whole date text with the literal\n// placeholder \"Date\", collapsing it into a single SigDate-styled run.\nconst body = context.document.body;\nconst results = body.search(\"1st October 2024\", { matchCase: true, matchWholeWord: false });\nresults.load(\"items\");\nawait context.sync();\n\nif (results.items.length > 0) {\n  results.items[0].insertText(\"Date\", \"Replace\");\n} else {\n  // Fallback: locate the SigBlock paragraph that still holds the old date\n  // text and replace its contents directly.\n  const paragraphs = body.paragraphs;\n  paragraphs.load(\"items/text\");\n  await context.sync();\n  for (const p of paragraphs.items) {\n    if (p.text.indexOf(\"October 2024\") !== -1) {\n      p.insertText(\"Date\", \"Replace\");\n      break;\n    }\n  }\n}\n\nawait context.sync();\n", "ps1": "# The signature block's date line currently reads \"1st October 2024\" and is\n# split across three runs (\"1\", \"st\", \" October 2024\") that all share the\n# SigDate character style. Replace that whole date text with the literal\n# placeholder \"Date\", collapsing it into a single SigDate-styled run.\n$d = $word.ActiveDocument\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Text = \"1st October 2024\"\n$find.Replacement.ClearFormatting()\n$find.Replacement.Text = \"Date\"\n$find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2)\n"}
